$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.90021176114846
$ws.Range("C2").Value = 5.330209548440743
$ws.Range("D2").Value = 6.020759350076236
$ws.Range("E2").Value = 10.43793165818004
$ws.Range("G2").Value = 46.39291866735508
$ws.Range("H2").Value = 18.5860839782774
$ws.Range("I2").Value = 28.85175970459092
$ws.Range("K2").Value = 12.5918413512403
$ws.Range("L2").Value = 10.17318235093565
$ws.Range("N2").Value = 21.23229170580766
$ws.Range("B3").Value = 14.68705456091702
$ws.Range("C3").Value = 5.145630984522278
$ws.Range("D3").Value = 5.908435913015204
$ws.Range("E3").Value = 10.44629139534924
$ws.Range("G3").Value = 46.29666572455955
$ws.Range("H3").Value = 18.61565104584932
$ws.Range("I3").Value = 28.8952854522726
$ws.Range("K3").Value = 12.44222970845614
$ws.Range("L3").Value = 10.16225546188473
$ws.Range("N3").Value = 21.29187332479179
$ws.Range("B4").Value = 14.55844593749071
$ws.Range("C4").Value = 5.027636119279143
$ws.Range("D4").Value = 5.840204445846966
$ws.Range("E4").Value = 10.45303375933569
$ws.Range("G4").Value = 46.24964604369455
$ws.Range("H4").Value = 18.63695230465499
$ws.Range("I4").Value = 28.92734068874117
$ws.Range("K4").Value = 12.35259927321732
$ws.Range("L4").Value = 10.15743497825742
$ws.Range("N4").Value = 21.3303403369765
$ws.Range("B5").Value = 14.50667575058209
$ws.Range("C5").Value = 4.978416087742681
$ws.Range("D5").Value = 5.812626119856652
$ws.Range("E5").Value = 10.45618635211789
$ws.Range("G5").Value = 46.23353197494709
$ws.Range("H5").Value = 18.64642312019163
$ws.Range("I5").Value = 28.94174155977483
$ws.Range("K5").Value = 12.31667825161544
$ws.Range("L5").Value = 10.15594726843372
$ws.Range("N5").Value = 21.34649049980858
$ws.Range("B6").Value = 14.49811990006581
$ws.Range("C6").Value = 4.970175662383449
$ws.Range("D6").Value = 5.808061659436435
$ws.Range("E6").Value = 10.45673430759657
$ws.Range("G6").Value = 46.23104047527448
$ws.Range("H6").Value = 18.64804344861745
$ws.Range("I6").Value = 28.94421356073894
$ws.Range("K6").Value = 12.31075131074335
$ws.Range("L6").Value = 10.15572907085471
$ws.Range("N6").Value = 21.34920090940291
$ws.Range("B7").Value = 14.55774507097078
$ws.Range("C7").Value = 5.026976869879751
$ws.Range("D7").Value = 5.839831542587238
$ws.Range("E7").Value = 10.45307463602322
$ws.Range("G7").Value = 46.24941637746689
$ws.Range("H7").Value = 18.63707683255547
$ws.Range("I7").Value = 28.92752948921245
$ws.Range("K7").Value = 12.35211232793686
$ws.Range("L7").Value = 10.15741298234133
$ws.Range("N7").Value = 21.33055622116502
$ws.Range("B8").Value = 14.82628485143868
$ws.Range("C8").Value = 5.267555163889893
$ws.Range("D8").Value = 5.981902455546345
$ws.Range("E8").Value = 10.44048026874133
$ws.Range("G8").Value = 46.35722989982846
$ws.Range("H8").Value = 18.59562501605034
$ws.Range("I8").Value = 28.8656598236049
$ws.Range("K8").Value = 12.53981816149846
$ws.Range("L8").Value = 10.169023904906
$ws.Range("N8").Value = 21.2524450271058
$ws.Range("B9").Value = 15.3678925899235
$ws.Range("C9").Value = 5.700941413780372
$ws.Range("D9").Value = 6.264528317746003
$ws.Range("E9").Value = 10.42853645491163
$ws.Range("G9").Value = 46.66395947127526
$ws.Range("H9").Value = 18.53934667708307
$ws.Range("I9").Value = 28.78671760684734
$ws.Range("K9").Value = 12.92368090084018
$ws.Range("L9").Value = 10.20668980254777
$ws.Range("N9").Value = 21.11417424279068
$ws.Range("B10").Value = 15.77091535306504
$ws.Range("C10").Value = 5.994356445431662
$ws.Range("D10").Value = 6.472253523150679
$ws.Range("E10").Value = 10.42750995356136
$ws.Range("G10").Value = 46.94649226726382
$ws.Range("H10").Value = 18.51329253261621
$ws.Range("I10").Value = 28.75466791990288
$ws.Range("K10").Value = 13.21268323190867
$ws.Range("L10").Value = 10.24330831033181
$ws.Range("N10").Value = 21.0216149422966
$ws.Range("B11").Value = 15.95452498709916
$ws.Range("C11").Value = 6.122143092691299
$ws.Range("D11").Value = 6.56629100660181
$ws.Range("E11").Value = 10.42871785559049
$ws.Range("G11").Value = 47.08717842356408
$ws.Range("H11").Value = 18.50476809295085
$ws.Range("I11").Value = 28.74573980073621
$ws.Range("K11").Value = 13.34510951031726
$ws.Range("L11").Value = 10.26187421158815
$ws.Range("N11").Value = 20.98145551138942
$ws.Range("B12").Value = 16.02402383395482
$ws.Range("C12").Value = 6.169695888264929
$ws.Range("D12").Value = 6.601796925143965
$ws.Range("E12").Value = 10.42941524335505
$ws.Range("G12").Value = 47.14217435524659
$ws.Range("H12").Value = 18.50201893980633
$ws.Range("I12").Value = 28.74317242282136
$ws.Range("K12").Value = 13.39534689192082
$ws.Range("L12").Value = 10.26917550166567
$ws.Range("N12").Value = 20.96652719984706
$ws.Range("B13").Value = 16.00905843618747
$ws.Range("C13").Value = 6.159492055438616
$ws.Range("D13").Value = 6.594155305681139
$ws.Range("E13").Value = 10.42925439009757
$ws.Range("G13").Value = 47.13025391845322
$ws.Range("H13").Value = 18.5025897176261
$ws.Range("I13").Value = 28.74368916309535
$ws.Range("K13").Value = 13.38452409607557
$ws.Range("L13").Value = 10.26759105498748
$ws.Range("N13").Value = 20.96972987218502
$ws.Range("B14").Value = 15.96024369289845
$ws.Range("C14").Value = 6.126072169304442
$ws.Range("D14").Value = 6.569214397099731
$ws.Range("E14").Value = 10.42877042584896
$ws.Range("G14").Value = 47.09166864581209
$ws.Range("H14").Value = 18.50453231986598
$ws.Range("I14").Value = 28.74551227359011
$ws.Range("K14").Value = 13.34924102586973
$ws.Range("L14").Value = 10.2624694876859
$ws.Range("N14").Value = 20.9802217585668
$ws.Range("B15").Value = 15.93033730541191
$ws.Range("C15").Value = 6.105491969488639
$ws.Range("D15").Value = 6.553922719873776
$ws.Range("E15").Value = 10.42850520895941
$ws.Range("G15").Value = 47.06825732838525
$ws.Range("H15").Value = 18.50578458990056
$ws.Range("I15").Value = 28.74673494280577
$ws.Range("K15").Value = 13.32763951176596
$ws.Range("L15").Value = 10.25936753303334
$ws.Range("N15").Value = 20.98668468141041
$ws.Range("B16").Value = 15.75891554460304
$ws.Range("C16").Value = 5.98588910814698
$ws.Range("D16").Value = 6.46609543260693
$ws.Range("E16").Value = 10.42746464573697
$ws.Range("G16").Value = 46.93754051042844
$ws.Range("H16").Value = 18.51391662026736
$ws.Range("I16").Value = 28.75536519454244
$ws.Range("K16").Value = 13.20404407549149
$ws.Range("L16").Value = 10.24213307910482
$ws.Range("N16").Value = 21.02427856098301
$ws.Range("B17").Value = 15.65377371385966
$ws.Range("C17").Value = 5.911044989131022
$ws.Range("D17").Value = 6.412071225940257
$ws.Range("E17").Value = 10.42725468881389
$ws.Range("G17").Value = 46.86044592781927
$ws.Range("H17").Value = 18.51975796395347
$ws.Range("I17").Value = 28.76210775184713
$ws.Range("K17").Value = 13.12843370113272
$ws.Range("L17").Value = 10.23204652380234
$ws.Range("N17").Value = 21.04783918449838
$ws.Range("B18").Value = 15.59332860384955
$ws.Range("C18").Value = 5.867462087435651
$ws.Range("D18").Value = 6.380957019360077
$ws.Range("E18").Value = 10.42729163179886
$ws.Range("G18").Value = 46.81724961050958
$ws.Range("H18").Value = 18.52343096072927
$ws.Range("I18").Value = 28.76651780514901
$ws.Range("K18").Value = 13.08503754805409
$ws.Range("L18").Value = 10.22642482307264
$ws.Range("N18").Value = 21.06157386216386
$ws.Range("B19").Value = 15.5728701719804
$ws.Range("C19").Value = 5.852614484149032
$ws.Range("D19").Value = 6.37041645143618
$ws.Range("E19").Value = 10.42733125040353
$ws.Range("G19").Value = 46.80282177441565
$ws.Range("H19").Value = 18.52472835537506
$ws.Range("I19").Value = 28.76810229873866
$ws.Range("K19").Value = 13.07036177001618
$ws.Range("L19").Value = 10.22455240000225
$ws.Range("N19").Value = 21.06625567896188
$ws.Range("B20").Value = 15.66496367693591
$ws.Range("C20").Value = 5.919067743657164
$ws.Range("D20").Value = 6.417826711894457
$ws.Range("E20").Value = 10.42726072313935
$ws.Range("G20").Value = 46.86853431488166
$ws.Range("H20").Value = 18.51910372379328
$ws.Range("I20").Value = 28.76133493691631
$ws.Range("K20").Value = 13.13647328100633
$ws.Range("L20").Value = 10.23310166919699
$ws.Range("N20").Value = 21.0453121572005
$ws.Range("B21").Value = 15.9745831265924
$ws.Range("C21").Value = 6.135911261484555
$ws.Range("D21").Value = 6.57654326522302
$ws.Range("E21").Value = 10.42890607227542
$ws.Range("G21").Value = 47.10295559831499
$ws.Range("H21").Value = 18.50394873224474
$ws.Range("I21").Value = 28.74495469858835
$ws.Range("K21").Value = 13.35960243788269
$ws.Range("L21").Value = 10.26396649704802
$ws.Range("N21").Value = 20.97713246526909
$ws.Range("B22").Value = 16.17673103038903
$ws.Range("C22").Value = 6.272743942367506
$ws.Range("D22").Value = 6.679650941577943
$ws.Range("E22").Value = 10.43137971979507
$ws.Range("G22").Value = 47.26618077994058
$ws.Range("H22").Value = 18.49683531687207
$ws.Range("I22").Value = 28.73899114374142
$ws.Range("K22").Value = 13.50593595644707
$ws.Range("L22").Value = 10.28571493095819
$ws.Range("N22").Value = 20.93420014850861
$ws.Range("B23").Value = 16.06888212735618
$ws.Range("C23").Value = 6.200166576966932
$ws.Range("D23").Value = 6.624689548633773
$ws.Range("E23").Value = 10.42993185514413
$ws.Range("G23").Value = 47.17815768107537
$ws.Range("H23").Value = 18.50037640985767
$ws.Range("I23").Value = 28.74173993881595
$ws.Range("K23").Value = 13.42780414472091
$ws.Range("L23").Value = 10.27396439096417
$ws.Range("N23").Value = 20.95696527669654
$ws.Range("B24").Value = 15.65990468204203
$ws.Range("C24").Value = 5.915442380583046
$ws.Range("D24").Value = 6.415224825016544
$ws.Range("E24").Value = 10.42725750390296
$ws.Range("G24").Value = 46.86487404430477
$ws.Range("H24").Value = 18.51939852537264
$ws.Range("I24").Value = 28.76168266419571
$ws.Range("K24").Value = 13.13283835512118
$ws.Range("L24").Value = 10.23262408575715
$ws.Range("N24").Value = 21.04645403628593
$ws.Range("B25").Value = 15.22018648757077
$ws.Range("C25").Value = 5.587986557767452
$ws.Range("D25").Value = 6.187888830422147
$ws.Range("E25").Value = 10.43040450785229
$ws.Range("G25").Value = 46.57086057952654
$ws.Range("H25").Value = 18.55188862188079
$ws.Range("I25").Value = 28.80352338119646
$ws.Range("K25").Value = 12.81841837655858
$ws.Range("L25").Value = 10.19491820424021
$ws.Range("N25").Value = 21.14999056211266
